$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: header values change
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: shift values - B2 gets new value, C2 deleted, D2/E2 updated
$ws.Range("B2").Value = 10.728384784212942
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 8.1335694119915729
$ws.Range("E2").Value = 13.937656838577013

# Row 3: B3 deleted, C3/E3 updated, D3 new value added
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 15.46825337714014
$ws.Range("D3").Value = 14.976920393192682
$ws.Range("E3").Value = 15.312650808290284

# Update the selection to match new used range B1:E3
$ws.Range("B1:E3").Select()
